$d = $word.ActiveDocument

# Primary footer (footer2.xml / rId13) holds the version & last-update text
$footer = $d.Sections.Item(1).Footers.Item(1)

# Update version number "4.0" -> "4.1"
$footer.Range.Find.Execute("Version 4.0", $true, $false, $false, $false, $false,
                            $true, 1, $false, "Version 4.1", 2)

# Update the last-update date field result text "2024-07-02" -> "2024-09-18"
$footer.Range.Find.Execute("2024-07-02", $true, $false, $false, $false, $false,
                            $true, 1, $false, "2024-09-18", 2)
